# Apply updated cryptocurrency data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new values are purely numeric-looking strings (e.g. "0.999")
# must be forced to remain text, otherwise Excel auto-converts them to numbers.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D20",
    "D21",
    "D22",
    "D23",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D35",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D44",
    "D45",
    "D47",
    "D48",
    "D50",
    "D51"
)

foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '67.596.94'
$ws.Range('E2').Value = '  +7.09%  '
$ws.Range('D3').Value = '3.546.82'
$ws.Range('E3').Value = '  +10.09%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '190.23'
$ws.Range('E5').Value = '  +9.36%  '
$ws.Range('D6').Value = '553.50'
$ws.Range('E6').Value = '  +3.91%  '
$ws.Range('D7').Value = '3.541.15'
$ws.Range('E7').Value = '  +10.01%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').Value = '0.635'
$ws.Range('E10').Value = '  +4.23%  '
$ws.Range('D11').Value = '0.156'
$ws.Range('E11').Value = '  +16.21%  '
$ws.Range('D12').Value = '55.05'
$ws.Range('E12').Value = '  +2.55%  '
$ws.Range('D13').Value = '0.0000272'
$ws.Range('E13').Value = '  +6.98%  '
$ws.Range('D14').Value = '9.41'
$ws.Range('E14').Value = '  +3.26%  '
$ws.Range('D15').Value = '4.098.68'
$ws.Range('E15').Value = '  +9.74%  '
$ws.Range('D16').Value = '3.536.79'
$ws.Range('E16').Value = '  +9.84%  '
$ws.Range('E17').Value = '  +4.31%  '
$ws.Range('D18').Value = '67.451.03'
$ws.Range('E18').Value = '  +7.25%  '
$ws.Range('E19').Value = '  +5.50%  '
$ws.Range('D20').Value = '11.94'
$ws.Range('E20').Value = '  +7.37%  '
$ws.Range('D21').Value = '0.994'
$ws.Range('D22').Value = '434.15'
$ws.Range('E22').Value = '  +18.04%  '
$ws.Range('D23').Value = '85.50'
$ws.Range('E23').Value = '  +5.18%  '
$ws.Range('E24').Value = '  +3.20%  '
$ws.Range('D25').Value = '4.15'
$ws.Range('E25').Value = '  +5.76%  '
$ws.Range('D26').Value = '11.13'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '2.91'
$ws.Range('E27').Value = '  +9.17%  '
$ws.Range('D28').Value = '12.11'
$ws.Range('E28').Value = '  +6.60%  '
$ws.Range('D29').Value = '9.02'
$ws.Range('E29').Value = '  +9.95%  '
$ws.Range('D30').Value = '30.40'
$ws.Range('E30').Value = '  +6.38%  '
$ws.Range('D31').Value = '644.93'
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('E32').Value = '  +3.00%  '
$ws.Range('D33').Value = '11.77'
$ws.Range('E33').Value = '  +3.90%  '
$ws.Range('E34').Value = '  +4.41%  '
$ws.Range('D35').Value = '59.85'
$ws.Range('E35').Value = '  +5.19%  '
$ws.Range('D36').Value = '0.0₃0834'
$ws.Range('E36').Value = '  +15.21%  '
$ws.Range('E37').Value = '  +4.73%  '
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').Value = '0.145'
$ws.Range('E39').Value = '  +16.67%  '
$ws.Range('D40').Value = '0.392'
$ws.Range('E40').Value = '  +3.77%  '
$ws.Range('D41').Value = '3.36'
$ws.Range('E41').Value = '  +13.78%  '
$ws.Range('D42').Value = '0.997'
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').Value = '3.043.67'
$ws.Range('E43').Value = '  +5.18%  '
$ws.Range('D44').Value = '2.66'
$ws.Range('E44').Value = '  +3.90%  '
$ws.Range('D45').Value = '2.89'
$ws.Range('E45').Value = '  +10.82%  '
$ws.Range('E46').Value = '  +6.99%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0419'
$ws.Range('E47').Value = '  +6.06%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '3.29'
$ws.Range('E48').Value = '  +9.50%  '
$ws.Range('E49').Value = '  +5.18%  '
$ws.Range('D50').Value = '8.75'
$ws.Range('E50').Value = '  +12.54%  '
$ws.Range('D51').Value = '140.93'
$ws.Range('E51').Value = '  +3.96%  '

# Restore default style on the forced-text cells so no stray cell style/format
# is left behind (matches original workbook formatting, which left these
# cells unstyled).
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
